$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) Merge the two runs that spell out "test_send_message - Test if
#    messages can be sent and received" into a single run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "test_send_message - Test if messages can be sent and received",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "test_send_message - Test if messages can be sent and received", 2) | Out-Null

# -----------------------------------------------------------------------
# 2) Merge the two runs of the overflow-wait bullet into one run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "test_message_overflow_wait - Tests if programs that chose to wait until able to send a message behave properly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "test_message_overflow_wait - Tests if programs that chose to wait until able to send a message behave properly", 2) | Out-Null

# -----------------------------------------------------------------------
# 3) Merge the two runs of the empty-mailbox bullet into one run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "test_recieve_empty_mailbox - Tests getting message from empty mailbox",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "test_recieve_empty_mailbox - Tests getting message from empty mailbox", 2) | Out-Null

# -----------------------------------------------------------------------
# 4) Drop the "_GoBack" bookmark that currently sits in the middle of the
#    "mailbox_errror_test:" heading (it is re-created at the end of the
#    list further down).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "mailbox_errror_test: Try to invoke all the possible ways to get an error",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mailbox_errror_test: Try to invoke all the possible ways to get an error", 2) | Out-Null

# -----------------------------------------------------------------------
# 5) Merge the many runs of the "rapid_fire_send_and_throw_an_exit_in_there"
#    bullet into a single trailing run.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "rapid_fire_send_and_throw_an_exit_in_there - Hopefully this can invoke the Mailbox dereference race condition. This invokes the pointer dereference race condition every once in a while, not the best test...a better test would fork and rerun this several times. If this test fails, you get a kernel oops.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "rapid_fire_send_and_throw_an_exit_in_there - Hopefully this can invoke the Mailbox dereference race condition. This invokes the pointer dereference race condition every once in a while, not the best test...a better test would fork and rerun this several times. If this test fails, you get a kernel oops.", 2) | Out-Null

# -----------------------------------------------------------------------
# 6) Merge the trailing runs of the
#    "rapid_fire_send_recieve_track_how_many_messages_we_get_eventaully"
#    bullet (keeping the inner "receiving" run separate, exactly as the
#    target keeps it as its own run).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    " - Stress test sending and ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " - Stress test sending and ", 2) | Out-Null

$d.Content.Find.Execute(
    " messages, one way. Hypothetically, two way is the same thing, just need to create two threads in each process",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " messages, one way. Hypothetically, two way is the same thing, just need to create two threads in each process", 2) | Out-Null

# -----------------------------------------------------------------------
# 7) Add the new bullet about "the_crazy_test_that_is_suggested_in_the_pdf_handout"
#    right after the "rapid_fire_send_recieve_track_..." bullet, re-using
#    the same ListParagraph/numId-10 formatting, and move the "_GoBack"
#    bookmark down to the end of this new paragraph.
# -----------------------------------------------------------------------
$i = 1
$targetIndex = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("rapid_fire_send_recieve_track_how_many_messages_we_get_eventaully")) {
        $targetIndex = $i
    }
    $i = $i + 1
}

$srcPara = $d.Paragraphs.Item($targetIndex)
$srcPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "the_crazy_test_that_is_suggested_in_the_pdf_handout"
$newPara.Range.InsertAfter(" " + [char]0x2013 + " stress test in the pdf. Made in like 20 minutes, doesn" + [char]0x2019 + "t clean up threads properly, need to sigint once cpu usage drops to 0")

$bookmarkPos = $newPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos)) | Out-Null

# -----------------------------------------------------------------------
# 8) Mark a rendered page break right before the "./" run that precedes
#    "mailbox_error_test" in the Output section.
# -----------------------------------------------------------------------
$i = 1
$outputIndex = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "./mailbox_error_test`r") {
        $outputIndex = $i
    }
    $i = $i + 1
}
$outPara = $d.Paragraphs.Item($outputIndex)
$outPara.Range.Characters.First.InsertBefore([char]11)
